$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.315.06'
$ws.Range("E2").Value = '  -2.34%  '

$ws.Range("D3").Value = '1.779.01'
$ws.Range("E3").Value = '  -0.72%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.37%  '

$ws.Range("B5").Value = 'USDC'
$ws.Range("C5").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.003'
$ws.Range("E5").Value = '  -0.10%  '

$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '305.74'
$ws.Range("E6").Value = '  -0.62%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4227'
$ws.Range("E7").Value = '  +1.46%  '

$ws.Range("E8").Value = '  +1.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07083'
$ws.Range("E9").Value = '  +0.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8383'
$ws.Range("E10").Value = '  -0.24%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.36'
$ws.Range("E11").Value = '  +1.41%  '

$ws.Range("D12").Value = '1.741.28'
$ws.Range("E12").Value = '  -6.95%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.445'
$ws.Range("E13").Value = '  +1.93%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.213'
$ws.Range("E14").Value = '  -1.27%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06903'
$ws.Range("E15").Value = '  +2.49%  '

$ws.Range("E16").Value = '  +0.03%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '79.00'
$ws.Range("E17").Value = '  -0.52%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008688'
$ws.Range("E18").Value = '  +0.08%  '

$ws.Range("E19").Value = '  -0.34%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.84'
$ws.Range("E20").Value = '  -1.10%  '

$ws.Range("D21").Value = '26.298.01'
$ws.Range("E21").Value = '  -3.37%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.112'
$ws.Range("E22").Value = '  +1.53%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.96'
$ws.Range("E23").Value = '  -0.40%  '

$ws.Range("D24").Value = '1.959.26'
$ws.Range("E24").Value = '  -6.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.32'
$ws.Range("E25").Value = '  -0.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.801'
$ws.Range("E26").Value = '  -7.05%  '

$ws.Range("E27").Value = '  -0.17%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.070'
$ws.Range("E28").Value = '  +1.72%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.18'
$ws.Range("E29").Value = '  +1.38%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.832'
$ws.Range("E30").Value = '  +12.56%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08868'
$ws.Range("E31").Value = '  -0.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7250'
$ws.Range("E32").Value = '  +1.89%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.126'
$ws.Range("E33").Value = '  +5.57%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.301'
$ws.Range("E34").Value = '  +0.34%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.749'
$ws.Range("E35").Value = '  -3.29%  '

$ws.Range("E36").Value = '  -0.20%  '

$ws.Range("E37").Value = '  +3.09%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05109'
$ws.Range("E38").Value = '  +0.43%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01887'
$ws.Range("E39").Value = '  -0.39%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.1610'
$ws.Range("E40").Value = '  -0.46%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4906'
$ws.Range("E41").Value = '  -0.37%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.592'
$ws.Range("E42").Value = '  +0.76%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.320'
$ws.Range("E43").Value = '  +5.36%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.087'
$ws.Range("E44").Value = '  +1.32%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.27'
$ws.Range("E45").Value = '  +0.71%  '

$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '104.60'
$ws.Range("E46").Value = '  +0.52%  '

$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.003'
$ws.Range("E47").Value = '  -0.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06167'
$ws.Range("E48").Value = '  -2.05%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.612'
$ws.Range("E49").Value = '  +1.67%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4453'
$ws.Range("E50").Value = '  -0.86%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.736'
$ws.Range("E51").Value = '  +4.48%  '
